# Applies a cyclic re-shuffle of the per-record columns (Fecha, Calidad,
# Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Unidad de comercializacion, Precio $/Kg, Kg / unidad) across rows 2-10
# of the active sheet. Columns A,B,C,E,F,G,H,I,J,K,R are identical across
# these rows and therefore remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to be written into each row (r2..r10), derived from the
# source rows of the pre-edit workbook following the mapping:
#   new2<-old6, new3<-old10, new4<-old7, new5<-old8, new6<-old2,
#   new7<-old3, new8<-old9, new9<-old4, new10<-old5
$rows = @{
    2  = @{ D = 44309; L = "Primera"; M = 10; N = 1600;  O = 1600;  P = 1600;  Q = '$/kilo (en caja de 15 kilos)'; S = 1600; T = 1 }
    3  = @{ D = 44400; L = "Primera"; M = 25; N = 1500;  O = 1500;  P = 1500;  Q = '$/kilo (en caja de 15 kilos)'; S = 1500; T = 1 }
    4  = @{ D = 44371; L = "Primera"; M = 20; N = 1800;  O = 1800;  P = 1800;  Q = '$/kilo (en caja de 15 kilos)'; S = 1800; T = 1 }
    5  = @{ D = 44371; L = "Segunda"; M = 30; N = 1200;  O = 1200;  P = 1200;  Q = '$/kilo (en caja de 15 kilos)'; S = 1200; T = 1 }
    6  = @{ D = 44336; L = "Primera"; M = 10; N = 1500;  O = 1500;  P = 1500;  Q = '$/kilo (en caja de 15 kilos)'; S = 1500; T = 1 }
    7  = @{ D = 44195; L = "Primera"; M = 20; N = 15000; O = 15000; P = 15000; Q = '$/bandeja 10 kilos'; S = 1500; T = 10 }
    8  = @{ D = 44292; L = "Primera"; M = 50; N = 14000; O = 14000; P = 14000; Q = '$/bandeja 10 kilos'; S = 1400; T = 10 }
    9  = @{ D = 44391; L = "Primera"; M = 15; N = 1500;  O = 1500;  P = 1500;  Q = '$/kilo (en caja de 15 kilos)'; S = 1500; T = 1 }
    10 = @{ D = 44391; L = "Segunda"; M = 20; N = 1000;  O = 1000;  P = 1000;  Q = '$/kilo (en caja de 15 kilos)'; S = 1000; T = 1 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $vals.S   # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $vals.T   # T - Kg / unidad
}
